# Bug fix in Eduati data files: CAR1_noCTRL_meas.xlsx
#
# Sheet1 ("measurements") had 43 stray trailing rows (45:87) that only
# carried a leftover index number in column A and no real measurement
# data -- left over from an earlier version of the sheet. Delete them so
# the real data block (rows 1:44, same row count as Sheet2 / Sheet3) is
# all that remains.
#
# Also update the saved view/selection state: Sheet1 becomes the active
# sheet (it was Sheet3), with F61 selected; Sheet3 is no longer the
# selected tab (its own selection, A2:N44, is left as-is).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Drop the stray trailing rows from Sheet1.
$ws1.Rows("45:87").Delete() | Out-Null

# Make Sheet1 the active sheet/tab and select F61, matching the
# workbook's final saved view state.
$ws1.Activate() | Out-Null
$ws1.Range("F61").Select() | Out-Null
